$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a real run-boundary split at an (already inserted) piece of
# text that has inherited identical formatting from its neighbour, by
# toggling Bold on/off. The engine merges runs with identical rPr, so a
# genuine (even if momentary) property change is required to keep the new
# text in its own <w:r>; restoring Bold to its original value (0 = not
# bold, which is what every run touched below already is) keeps the visible
# formatting unchanged.
# ---------------------------------------------------------------------------

# ============================================================
# Edit 1: paragraph with the merged "samplecucumber" run.
# Split it into "sample" / "maven" / bookmark(_GoBack) / "cucumber".
# ============================================================
$p1 = $d.Paragraphs.Item(40)
$r1 = $p1.Range
$full1 = $r1.Text
$idx1 = $full1.IndexOf("samplecucumber")
$absStart1 = $r1.Start + $idx1

# Shrink "samplecucumber" down to just "sample" (keeps the original run).
$target1 = $d.Range($absStart1, $absStart1 + ([string]"samplecucumber").Length)
$target1.Find.Execute("samplecucumber", $true, $false, $false, $false, $false, `
                       $true, 0, $false, "sample", 2) | Out-Null

$sampleEnd1 = $absStart1 + 6
$mavenStart1 = $sampleEnd1
$mavenEnd1 = $mavenStart1 + 5
$cucStart1 = $mavenEnd1
$cucEnd1 = $cucStart1 + 8

# Insert "maven" + "cucumber" text together right after "sample" so the new
# text inherits "sample"'s character formatting.
$afterSample1 = $d.Range($sampleEnd1, $sampleEnd1)
$afterSample1.InsertAfter("mavencucumber")

# Split "maven" into its own run.
$mavenRange1 = $d.Range($mavenStart1, $mavenEnd1)
$mavenRange1.Font.Bold = 1
$mavenRange1.Font.Bold = 0

# Split "cucumber" into its own run.
$cucRange1 = $d.Range($cucStart1, $cucEnd1)
$cucRange1.Font.Bold = 1
$cucRange1.Font.Bold = 0

# Move the singleton "_GoBack" bookmark to sit between "maven" and "cucumber".
$bmRange1 = $d.Range($mavenEnd1, $mavenEnd1)
$d.Bookmarks.Add("_GoBack", $bmRange1) | Out-Null

# ============================================================
# Edit 2: paragraph that already has separate "sample" / "cucumber" runs.
# Insert a new "maven" run between them.
# ============================================================
$p2 = $d.Paragraphs.Item(33)
$r2 = $p2.Range
$full2 = $r2.Text
$idx2 = $full2.IndexOf("cucumber")
$cucStart2 = $r2.Start + $idx2

$insertPoint2 = $d.Range($cucStart2, $cucStart2)
$insertPoint2.InsertAfter("maven")

$mavenStart2 = $cucStart2
$mavenEnd2 = $mavenStart2 + 5
$mavenRange2 = $d.Range($mavenStart2, $mavenEnd2)
$mavenRange2.Font.Bold = 1
$mavenRange2.Font.Bold = 0

Write-Host "Done"
